$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '46.851.37'
$ws.Range('E2').Value = '  +6.86%  '
$ws.Range('D3').Value = '2.312.07'
$ws.Range('E3').Value = '  +5.68%  '
$ws.Range('E4').Value = '  -0.45%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '298.95'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.86%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.99'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +13.31%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.572'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.81%  '
$ws.Range('E8').Value = '  -0.39%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.528'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +10.81%  '
$ws.Range('E10').Value = '  +11.87%  '
$ws.Range('E11').Value = '  +5.09%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.32'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +9.40%  '
$ws.Range('E13').Value = '  +2.16%  '
$ws.Range('D14').Value = '2.664.41'
$ws.Range('E14').Value = '  +5.80%  '
$ws.Range('D15').Value = '2.311.75'
$ws.Range('E15').Value = '  +2.38%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.02'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +9.09%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.817'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +7.57%  '
$ws.Range('D18').Value = '46.787.08'
$ws.Range('E18').Value = '  +8.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.17'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +23.79%  '
$ws.Range('E20').Value = '  +7.77%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.14'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +6.18%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '66.77'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +7.16%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '248.04'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +8.99%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.92'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +6.27%  '
$ws.Range('E25').Value = '  +10.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.21%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '42.81'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +21.35%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.24'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.87'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +8.64%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '20.17'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +6.34%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.74'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +9.88%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '147.09'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0799'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +10.89%  '
$ws.Range('E34').Value = '  +5.21%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.11'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +8.63%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.113'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +12.91%  '
$ws.Range('E37').Value = '  +3.23%  '
$ws.Range('E38').Value = '  +10.66%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '15.63'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +19.40%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.02'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +16.06%  '
$ws.Range('E41').Value = '  +12.96%  '
$ws.Range('E42').Value = '  +10.44%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.999'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.42%  '
$ws.Range('E44').Value = '  +21.70%  '
$ws.Range('D45').Value = '1.842.67'
$ws.Range('E45').Value = '  +5.43%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '90.08'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +24.78%  '
$ws.Range('E47').Value = '  +17.83%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '76.19'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +16.69%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.95'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +11.47%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '97.17'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '54.21'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +12.41%  '
